$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated daily sales ("faturamento") rows: Dia, total_venda, Mes, Ano, Periodo
$rows = @(
  @(1,19371.37,4,2025,"04/2025"),
  @(2,29218.71,4,2025,"04/2025"),
  @(3,19819.3,4,2025,"04/2025"),
  @(4,25399.91,4,2025,"04/2025"),
  @(5,18509.5,4,2025,"04/2025"),
  @(6,5823.21,4,2025,"04/2025"),
  @(7,21764.19,4,2025,"04/2025"),
  @(8,38515.34,4,2025,"04/2025"),
  @(9,21733.46,4,2025,"04/2025"),
  @(10,25715.4,4,2025,"04/2025"),
  @(11,22833.25,4,2025,"04/2025"),
  @(12,21596.15,4,2025,"04/2025"),
  @(13,6451,4,2025,"04/2025"),
  @(14,24085.6,4,2025,"04/2025"),
  @(15,15644.59,4,2025,"04/2025"),
  @(16,18372.2,4,2025,"04/2025"),
  @(17,23446.19,4,2025,"04/2025"),
  @(18,15474.35,4,2025,"04/2025"),
  @(19,9164.8,4,2025,"04/2025"),
  @(20,2719,4,2025,"04/2025"),
  @(21,15666.86,4,2025,"04/2025"),
  @(22,15219.82,4,2025,"04/2025"),
  @(23,19895.36,4,2025,"04/2025"),
  @(24,19632.82,4,2025,"04/2025"),
  @(25,19212.07,4,2025,"04/2025"),
  @(26,14442,4,2025,"04/2025"),
  @(27,5289.3,4,2025,"04/2025"),
  @(28,16598.71,4,2025,"04/2025"),
  @(1,13803.81,3,2025,"03/2025"),
  @(2,4705.9,3,2025,"03/2025"),
  @(3,16680.45,3,2025,"03/2025"),
  @(4,7657.8,3,2025,"03/2025"),
  @(5,37984.79,3,2025,"03/2025"),
  @(6,19791.3,3,2025,"03/2025"),
  @(7,51663.12,3,2025,"03/2025"),
  @(8,14212.1,3,2025,"03/2025"),
  @(9,8836.05,3,2025,"03/2025"),
  @(10,38138.8,3,2025,"03/2025"),
  @(11,30745.57,3,2025,"03/2025"),
  @(12,44416.44,3,2025,"03/2025"),
  @(13,20045.79,3,2025,"03/2025"),
  @(14,23677.64,3,2025,"03/2025"),
  @(15,23636.62,3,2025,"03/2025"),
  @(16,4642,3,2025,"03/2025"),
  @(17,37009.74,3,2025,"03/2025"),
  @(18,15851.63,3,2025,"03/2025"),
  @(19,33313.16,3,2025,"03/2025"),
  @(20,23394.13,3,2025,"03/2025"),
  @(21,17765.24,3,2025,"03/2025"),
  @(22,21750.9,3,2025,"03/2025"),
  @(23,5923.91,3,2025,"03/2025"),
  @(24,23370.46,3,2025,"03/2025"),
  @(25,27227.58,3,2025,"03/2025"),
  @(26,20282.52,3,2025,"03/2025"),
  @(27,19695.38,3,2025,"03/2025"),
  @(28,17250.36,3,2025,"03/2025"),
  @(29,15973.05,3,2025,"03/2025"),
  @(30,6906.9,3,2025,"03/2025"),
  @(31,24267.47,3,2025,"03/2025"),
  @(1,22837.82,2,2025,"02/2025"),
  @(2,10035.5,2,2025,"02/2025"),
  @(3,23114.36,2,2025,"02/2025"),
  @(4,17250.02,2,2025,"02/2025"),
  @(5,28584.32,2,2025,"02/2025"),
  @(6,25749.58,2,2025,"02/2025"),
  @(7,54601.36,2,2025,"02/2025"),
  @(8,25472.3,2,2025,"02/2025"),
  @(9,6812.99,2,2025,"02/2025"),
  @(10,34211.97,2,2025,"02/2025"),
  @(11,20510.6,2,2025,"02/2025"),
  @(12,28660.12,2,2025,"02/2025"),
  @(13,35820.62,2,2025,"02/2025"),
  @(14,20762.42,2,2025,"02/2025"),
  @(15,18502.7,2,2025,"02/2025"),
  @(16,4694,2,2025,"02/2025"),
  @(17,57356.01,2,2025,"02/2025"),
  @(18,38766.64,2,2025,"02/2025"),
  @(19,37172.24,2,2025,"02/2025"),
  @(20,22352.55,2,2025,"02/2025"),
  @(21,38289.45,2,2025,"02/2025"),
  @(22,19482.2,2,2025,"02/2025"),
  @(23,8507.9,2,2025,"02/2025"),
  @(24,45768.1,2,2025,"02/2025"),
  @(25,41476.67,2,2025,"02/2025"),
  @(26,20446.88,2,2025,"02/2025"),
  @(27,29037.01,2,2025,"02/2025"),
  @(28,26934.44,2,2025,"02/2025"),
  @(2,21267.77,1,2025,"01/2025"),
  @(3,24209.45,1,2025,"01/2025"),
  @(4,14089.99,1,2025,"01/2025"),
  @(5,9388,1,2025,"01/2025"),
  @(6,31129.04,1,2025,"01/2025"),
  @(7,22066.68,1,2025,"01/2025"),
  @(8,22417.6,1,2025,"01/2025"),
  @(9,23072.95,1,2025,"01/2025"),
  @(10,17291.29,1,2025,"01/2025"),
  @(11,15776,1,2025,"01/2025"),
  @(12,4529.75,1,2025,"01/2025"),
  @(13,30004,1,2025,"01/2025"),
  @(14,16520.27,1,2025,"01/2025"),
  @(15,24823.94,1,2025,"01/2025"),
  @(16,31800.92,1,2025,"01/2025"),
  @(17,26415.47,1,2025,"01/2025"),
  @(18,20240.01,1,2025,"01/2025"),
  @(19,9165,1,2025,"01/2025"),
  @(20,20439.55,1,2025,"01/2025"),
  @(21,21938.12,1,2025,"01/2025"),
  @(22,17800.9,1,2025,"01/2025"),
  @(23,17359.35,1,2025,"01/2025"),
  @(24,21841.03,1,2025,"01/2025"),
  @(25,15414.51,1,2025,"01/2025"),
  @(26,6505.9,1,2025,"01/2025"),
  @(27,19682.44,1,2025,"01/2025"),
  @(28,22828.47,1,2025,"01/2025"),
  @(29,15542.37,1,2025,"01/2025"),
  @(30,22876.75,1,2025,"01/2025"),
  @(31,21144.84,1,2025,"01/2025")
)

$n = $rows.Length
$data = New-Object 'object[,]' $n,5
for ($i = 0; $i -lt $n; $i++) {
    for ($j = 0; $j -lt 5; $j++) {
        $data[$i, $j] = $rows[$i][$j]
    }
}

$startRow = 2
$endRow = $startRow + $n - 1
$rangeAddr = "A" + $startRow + ":E" + $endRow
$range = $ws.Range($rangeAddr)
$range.Value = $data
